# Update "想去人数" (want-to-go count) figures refreshed by the gh-pages
# data regeneration (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 7501
$wsExhibit.Range("F6").Value = 450
$wsExhibit.Range("F7").Value = 4101
$wsExhibit.Range("F8").Value = 323
$wsExhibit.Range("F9").Value = 576
$wsExhibit.Range("F11").Value = 656
$wsExhibit.Range("F12").Value = 141

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7501
$wsAll.Range("F8").Value = 450
$wsAll.Range("F9").Value = 4101
$wsAll.Range("F10").Value = 323
$wsAll.Range("F11").Value = 576
$wsAll.Range("F13").Value = 656
$wsAll.Range("F15").Value = 141
